$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MAXQDA 12")

# New coded segments exported from MAXQDA for document 16897 (coder: white)
# Append as rows 577-585, matching the formatting of the last existing row (576)

$srcRow = $ws.Range("A576:M576")
for ($i = 577; $i -le 585; $i++) {
    $dstRow = $ws.Range("A" + $i + ":M" + $i)
    $srcRow.Copy()
    $dstRow.PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = 0

# Row 577 - Patient:Age
$ws.Cells.Item(577, 1).Value = "●"
$ws.Cells.Item(577, 4).Value = "16897"
$ws.Cells.Item(577, 5).Value = "Patient:Age"
$ws.Cells.Item(577, 6).Value = "1: 2832"
$ws.Cells.Item(577, 7).Value = "1: 2833"
$ws.Cells.Item(577, 8).Value = 0
$ws.Cells.Item(577, 9).Value = "57"
$ws.Cells.Item(577, 10).Value = 2
$ws.Cells.Item(577, 11).Value = 0.013185654008438817
$ws.Cells.Item(577, 12).Value = "white"
$ws.Cells.Item(577, 13).Value = "08/21/2019 15:29:14"

# Row 578 - Patient:Sex
$ws.Cells.Item(578, 1).Value = "●"
$ws.Cells.Item(578, 4).Value = "16897"
$ws.Cells.Item(578, 5).Value = "Patient:Sex"
$ws.Cells.Item(578, 6).Value = "1: 2844"
$ws.Cells.Item(578, 7).Value = "1: 2847"
$ws.Cells.Item(578, 8).Value = 0
$ws.Cells.Item(578, 9).Value = "male"
$ws.Cells.Item(578, 10).Value = 4
$ws.Cells.Item(578, 11).Value = 0.026371308016877634
$ws.Cells.Item(578, 12).Value = "white"
$ws.Cells.Item(578, 13).Value = "08/21/2019 15:29:21"

# Row 579 - Patient:Comorbidities
$ws.Cells.Item(579, 1).Value = "●"
$ws.Cells.Item(579, 4).Value = "16897"
$ws.Cells.Item(579, 5).Value = "Patient:Comorbidities"
$ws.Cells.Item(579, 6).Value = "1: 2850"
$ws.Cells.Item(579, 7).Value = "1: 2857"
$ws.Cells.Item(579, 8).Value = 0
$ws.Cells.Item(579, 9).Value = "diabetic"
$ws.Cells.Item(579, 10).Value = 8
$ws.Cells.Item(579, 11).Value = 0.052742616033755269
$ws.Cells.Item(579, 12).Value = "white"
$ws.Cells.Item(579, 13).Value = "08/21/2019 15:29:29"

# Row 580 - Bacteria:Binomial (genus species)
$ws.Cells.Item(580, 1).Value = "●"
$ws.Cells.Item(580, 4).Value = "16897"
$ws.Cells.Item(580, 5).Value = "Bacteria:Binomial (genus species)"
$ws.Cells.Item(580, 6).Value = "1: 146"
$ws.Cells.Item(580, 7).Value = "1: 164"
$ws.Cells.Item(580, 8).Value = 0
$ws.Cells.Item(580, 9).Value = "Nocardia asteroides"
$ws.Cells.Item(580, 10).Value = 19
$ws.Cells.Item(580, 11).Value = 0.12526371308016879
$ws.Cells.Item(580, 12).Value = "white"
$ws.Cells.Item(580, 13).Value = "08/21/2019 15:30:29"

# Row 581 - Event month
$ws.Cells.Item(581, 1).Value = "●"
$ws.Cells.Item(581, 4).Value = "16897"
$ws.Cells.Item(581, 5).Value = "Event month"
$ws.Cells.Item(581, 6).Value = "1: 3552"
$ws.Cells.Item(581, 7).Value = "1: 3555"
$ws.Cells.Item(581, 8).Value = 0
$ws.Cells.Item(581, 9).Value = "July"
$ws.Cells.Item(581, 10).Value = 4
$ws.Cells.Item(581, 11).Value = 0.026371308016877634
$ws.Cells.Item(581, 12).Value = "white"
$ws.Cells.Item(581, 13).Value = "08/21/2019 15:30:48"

# Row 582 - Event year
$ws.Cells.Item(582, 1).Value = "●"
$ws.Cells.Item(582, 4).Value = "16897"
$ws.Cells.Item(582, 5).Value = "Event year"
$ws.Cells.Item(582, 6).Value = "1: 3557"
$ws.Cells.Item(582, 7).Value = "1: 3560"
$ws.Cells.Item(582, 8).Value = 0
$ws.Cells.Item(582, 9).Value = "2009"
$ws.Cells.Item(582, 10).Value = 4
$ws.Cells.Item(582, 11).Value = 0.026371308016877634
$ws.Cells.Item(582, 12).Value = "white"
$ws.Cells.Item(582, 13).Value = "08/21/2019 15:30:52"

# Row 583 - Drug Resisted
$ws.Cells.Item(583, 1).Value = "●"
$ws.Cells.Item(583, 4).Value = "16897"
$ws.Cells.Item(583, 5).Value = "Drug Resisted"
$ws.Cells.Item(583, 6).Value = "1: 5495"
$ws.Cells.Item(583, 7).Value = "1: 5523"
$ws.Cells.Item(583, 8).Value = 0
$ws.Cells.Item(583, 9).Value = "trimethoprim–sulfamethoxazole"
$ws.Cells.Item(583, 10).Value = 29
$ws.Cells.Item(583, 11).Value = 0.19119198312236288
$ws.Cells.Item(583, 12).Value = "white"
$ws.Cells.Item(583, 13).Value = "08/21/2019 15:31:56"

# Row 584 - Patient:Outcome
$ws.Cells.Item(584, 1).Value = "●"
$ws.Cells.Item(584, 4).Value = "16897"
$ws.Cells.Item(584, 5).Value = "Patient:Outcome"
$ws.Cells.Item(584, 6).Value = "1: 6298"
$ws.Cells.Item(584, 7).Value = "1: 6301"
$ws.Cells.Item(584, 8).Value = 0
$ws.Cells.Item(584, 9).Value = "died"
$ws.Cells.Item(584, 10).Value = 4
$ws.Cells.Item(584, 11).Value = 0.026371308016877634
$ws.Cells.Item(584, 12).Value = "white"
$ws.Cells.Item(584, 13).Value = "08/21/2019 15:33:26"

# Row 585 - Location:Country
$ws.Cells.Item(585, 1).Value = "●"
$ws.Cells.Item(585, 4).Value = "16897"
$ws.Cells.Item(585, 5).Value = "Location:Country"
$ws.Cells.Item(585, 6).Value = "1: 3235"
$ws.Cells.Item(585, 7).Value = "1: 3239"
$ws.Cells.Item(585, 8).Value = 0
$ws.Cells.Item(585, 9).Value = "India"
$ws.Cells.Item(585, 10).Value = 5
$ws.Cells.Item(585, 11).Value = 0.032964135021097046
$ws.Cells.Item(585, 12).Value = "white"
$ws.Cells.Item(585, 13).Value = "08/21/2019 15:36:26"
